$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: "FWPKG" divider row, same box-border style as the rest of the table (style index 3) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A25:C25").PasteSpecial(-4122) | Out-Null

$ws.Range("A25").Value = "FWPKG"
$ws.Range("B25").Value = "FWPKG"
$ws.Range("C25").Value = "FWPKG"

# --- Rows 36-58: duplicate of rows 2-24, all with the plain box-border style (style index 3) ---
$ws.Range("A36:C58").PasteSpecial(-4122) | Out-Null

$tableData = @(
    @("Mellanox InfiniBand and Ethernet Driver", "Mellanox OFED", "MOFED"),
    @("Mellanox InfiniBand and Ethernet Driver for Microsoft Windows", "WinOF", "WINOF"),
    @("Mellanox WinOF-2 InfiniBand and Ethernet driver", "WinOF2", "WINOF2"),
    @("nmlx4_en Driver Component for VMware 6.7", "nmlx4_en ESXi 6.7 Driver Smart Components", "6.7_nmlx4_SC"),
    @("nmlx4_en Driver Component for VMware 6.5", "nmlx4_en ESXi 6.5 Driver Smart Components", "6.5_nmlx4_SC"),
    @("nmlx4_en Driver Component for VMware 6.0", "nmlx4_en ESXi 6.0 Driver Smart Components", "6.0_nmlx4_SC"),
    @("nmlx5_en Driver Component for VMware 7.0", "nmlx5_en ESXi 7.0 Driver Smart Components", "7.0_nmlx5_SC"),
    @("nmlx5_en Driver Component for VMware 6.7", "nmlx5_en ESXi 6.7 Driver Smart Components", "6.7_nmlx5_SC"),
    @("nmlx5_en Driver Component for VMware 6.5", "nmlx5_en ESXi 6.5 Driver Smart Components", "6.5_nmlx5_SC"),
    @("nmlx5_en Driver Component for VMware 6.0", "nmlx5_en ESXi 6.0 Driver Smart Components", "6.0_nmlx5_SC"),
    @("HPE Mellanox RoCE", "Linux RoCE driver Smart Components", "RoCE_SC"),
    @("Mellanox net-mst Kernel Driver Component for VMware ESXi 6.0", "net-mst Kernel Driver ESXi6.0 Smart component", "6.0_MST_SC"),
    @("Mellanox net-mst Kernel Driver Component for VMware ESXi 6.5 and 6.7", "net-mst Kernel Driver ESXi6.5 and 6.7 Smart component", "6.5_6.7_MST_SC"),
    @("Mellanox net-mst Kernel Driver Component for VMware ESXi 7.0", "net-mst Kernel Driver ESXi7.0 Smart component", "7.0_MST_SC"),
    @("HPE Mellanox MFT Driver and Firmware Tools", "Mellanox MFT DR_FW Tools Smart Components", "Linux_MFT_SC"),
    @("Online Firmware Upgrade Utility (Linux x86_64)", "Linux firmware smart components", "Linux_FW_SC"),
    @("Online Firmware Upgrade Utility (Windows x64)", "Windows firmware smart components", "Windows_FW_SC"),
    @("Online Firmware Upgrade Utility (ESXi 6.0) ", "ESXi 6.0 firmware smart components", "ESXi6.0_FW_SC"),
    @("Online Firmware Upgrade Utility (ESXi 6.5) ", "ESXi 6.5 firmware smart components", "ESXi6.5_FW_SC"),
    @("Online Firmware Upgrade Utility (ESXi 6.7) ", "ESXi 6.7 firmware smart components", "ESXi6.7_FW_SC"),
    @("Online Firmware Upgrade Utility (ESXi 7.0) ", "ESXi 7.0 firmware smart components", "ESXi7.0_FW_SC"),
    @("Firmware for ", "Firmware binary posting", "FW_Binary"),
    @("HPE Mellanox Firmware Tools", "Mellanox MFT  ", "MFT")
)

$startRow = 36
for ($i = 0; $i -lt $tableData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $tableData[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
}

# --- Row 59: "FWPKG" footer row with distinct yellow-fill styling ---
# A59: yellow fill + thin left/right border (new style)
$ws.Range("A59").Interior.Color = 65535
$ws.Range("A59").Borders.Item(7).LineStyle = 1
$ws.Range("A59").Borders.Item(10).LineStyle = 1
$ws.Range("A59").Value = "FWPKG"

# B59/C59: yellow fill, no border (new style)
$ws.Range("B59:C59").Interior.Color = 65535
$ws.Range("B59").Value = "FWPKG"
$ws.Range("C59").Value = "FWPKG"

# --- Update the active selection to match the authored state ---
$ws.Range("B8").Select() | Out-Null

Write-Host "Edit applied"
